$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain text, matching the original inline-string cells,
# so Excel does not auto-convert numeric-looking price strings into Number cells.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "24.213.91"
$ws.Range("D3").Value = "1.646.18"
$ws.Range("D4").Value = "0.9991"
$ws.Range("D5").Value = "309.59"
$ws.Range("D7").Value = "0.3912"
$ws.Range("D8").Value = "0.3867"
$ws.Range("D9").Value = "1.004"
$ws.Range("D10").Value = "1.360"
$ws.Range("D11").Value = "49.50"
$ws.Range("D12").Value = "0.08620"
$ws.Range("D13").Value = "23.57"
$ws.Range("D14").Value = "7.104"
$ws.Range("D15").Value = "0.00001292"
$ws.Range("D16").Value = "7.477"
$ws.Range("D17").Value = "1.640.87"
$ws.Range("D18").Value = "95.09"
$ws.Range("D19").Value = "0.06921"
$ws.Range("D20").Value = "20.37"
$ws.Range("D21").Value = "6.895"
$ws.Range("D22").Value = "1.001"
$ws.Range("D23").Value = "13.57"
$ws.Range("D24").Value = "24.200.60"
$ws.Range("D25").Value = "2.404"
$ws.Range("D26").Value = "2.837"
$ws.Range("D27").Value = "22.35"
$ws.Range("D28").Value = "157.85"
$ws.Range("D29").Value = "8.430"
$ws.Range("D30").Value = "5.371"
$ws.Range("D31").Value = "140.22"
$ws.Range("D32").Value = "2.406"
$ws.Range("D33").Value = "1.815.96"
$ws.Range("D34").Value = "6.956"
$ws.Range("D35").Value = "0.08108"
$ws.Range("D36").Value = "0.02906"
$ws.Range("D37").Value = "0.2688"
$ws.Range("D38").Value = "0.9500"
$ws.Range("D39").Value = "0.09203"
$ws.Range("D40").Value = "10.14"
$ws.Range("D41").Value = "1.459"
$ws.Range("D42").Value = "0.7558"
$ws.Range("D43").Value = "13.01"
$ws.Range("D44").Value = "15.98"
$ws.Range("D45").Value = "0.6902"
$ws.Range("D46").Value = "2.457"
$ws.Range("D47").Value = "4.095"
$ws.Range("D49").Value = "0.08377"
$ws.Range("D50").Value = "133.52"
$ws.Range("D51").Value = "1.258"

# Clear the temporary text format so the cells end up with no explicit style,
# matching the original workbook (which had no "s" attribute on these cells).
$ws.Range("D2:D51").ClearFormats()

# Volume(1h) column: values are already non-numeric-looking text (e.g. "  -2.27%  "),
# so a plain .Value assignment keeps them as text without any style side effects.
$ws.Range("E2").Value = "  -2.27%  "
$ws.Range("E3").Value = "  -2.14%  "
$ws.Range("E4").Value = "  -0.71%  "
$ws.Range("E5").Value = "  -1.28%  "
$ws.Range("E6").Value = "  -0.34%  "
$ws.Range("E7").Value = "  -0.29%  "
$ws.Range("E8").Value = "  -2.40%  "
$ws.Range("E9").Value = "  -0.19%  "
$ws.Range("E10").Value = "  -4.12%  "
$ws.Range("E11").Value = "  -4.19%  "
$ws.Range("E12").Value = "  -0.48%  "
$ws.Range("E13").Value = "  -6.08%  "
$ws.Range("E14").Value = "  -2.68%  "
$ws.Range("E15").Value = "  -1.86%  "
$ws.Range("E16").Value = "  -4.01%  "
$ws.Range("E17").Value = "  +4.69%  "
$ws.Range("E18").Value = "  +1.17%  "
$ws.Range("E19").Value = "  -2.85%  "
$ws.Range("E20").Value = "  +0.94%  "
$ws.Range("E21").Value = "  -3.16%  "
$ws.Range("E22").Value = "  -0.39%  "
$ws.Range("E23").Value = "  -3.60%  "
$ws.Range("E24").Value = "  -2.34%  "
$ws.Range("E25").Value = "  +2.25%  "
$ws.Range("E26").Value = "  +2.74%  "
$ws.Range("E27").Value = "  -5.23%  "
$ws.Range("E28").Value = "  -2.77%  "
$ws.Range("E29").Value = "  +8.24%  "
$ws.Range("E30").Value = "  -6.99%  "
$ws.Range("E31").Value = "  -6.40%  "
$ws.Range("E32").Value = "  -6.16%  "
$ws.Range("E33").Value = "  +4.24%  "
$ws.Range("E34").Value = "  +0.56%  "
$ws.Range("E35").Value = "  -3.63%  "
$ws.Range("E36").Value = "  -5.27%  "
$ws.Range("E37").Value = "  -3.95%  "
$ws.Range("E38").Value = "  -5.39%  "
$ws.Range("E39").Value = "  -3.58%  "
$ws.Range("E40").Value = "  -3.32%  "
$ws.Range("E41").Value = "  -0.75%  "
$ws.Range("E42").Value = "  -4.64%  "
$ws.Range("E43").Value = "  -4.66%  "
$ws.Range("E44").Value = "  -3.76%  "
$ws.Range("E45").Value = "  -3.47%  "
$ws.Range("E46").Value = "  -4.71%  "
$ws.Range("E47").Value = "  -2.13%  "
$ws.Range("E48").Value = "  -0.33%  "
$ws.Range("E49").Value = "  -3.90%  "
$ws.Range("E50").Value = "  -3.31%  "
$ws.Range("E51").Value = "  -5.75%  "
